# Generate Report for Handoff
# Updates the localization-status report for the bdc28621-...md file,
# reflecting that it is now "Ready for handoff" (rather than handed back).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the bdc28621-516e-4ab6-92f9-1cb06dee1b24.md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-24 04:46:14"

# --- zh-cn sheet: row 3 is the bdc28621-516e-4ab6-92f9-1cb06dee1b24.md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-24 04:46:08"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ff103b250cc315289c2edd3b2ed98142c539ca6/e2e/bdc28621-516e-4ab6-92f9-1cb06dee1b24.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c900d6eb20090dad32e317bb2f5dfb996c7cab2b/e2e/bdc28621-516e-4ab6-92f9-1cb06dee1b24.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is the bdc28621-516e-4ab6-92f9-1cb06dee1b24.md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-24 04:46:14"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ff103b250cc315289c2edd3b2ed98142c539ca6/e2e/bdc28621-516e-4ab6-92f9-1cb06dee1b24.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c900d6eb20090dad32e317bb2f5dfb996c7cab2b/e2e/bdc28621-516e-4ab6-92f9-1cb06dee1b24.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
